$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B13").Value = "WASD WASD"
$ws.Range("C13").Value = "ASD ASDASD"
